# Updated transition probability matrix values (West Virginia_B.xlsx)
# Reflects additional simulated games / re-run optimization logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1963470319634703
$ws.Range("C2").Value = 0.5342465753424658
$ws.Range("J2").Value = 0.0319634703196347
$ws.Range("P2").Value = 0.1689497716894977
$ws.Range("S2").Value = 0.0684931506849315
# Row 3
$ws.Range("B3").Value = 0.0078125
$ws.Range("C3").Value = 0.0078125
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.796875
$ws.Range("S3").Value = 0.15625
# Row 4
$ws.Range("J4").Value = 0.03448275862068965
$ws.Range("P4").Value = 0.7586206896551724
$ws.Range("S4").Value = 0.2068965517241379
# Row 6
$ws.Range("B6").Value = 0.06153846153846154
$ws.Range("F6").Value = 0.04615384615384616
$ws.Range("J6").Value = 0.2564102564102564
$ws.Range("O6").Value = 0.04102564102564103
$ws.Range("Q6").Value = 0.1128205128205128
$ws.Range("R6").Value = 0.09230769230769231
$ws.Range("S6").Value = 0.3897435897435897
# Row 7
$ws.Range("B7").Value = 0.08620689655172414
$ws.Range("D7").Value = 0.005747126436781609
$ws.Range("E7").Value = 0.005747126436781609
$ws.Range("F7").Value = 0.06321839080459771
$ws.Range("J7").Value = 0.1781609195402299
$ws.Range("O7").Value = 0.02298850574712644
$ws.Range("Q7").Value = 0.2068965517241379
$ws.Range("R7").Value = 0.08620689655172414
$ws.Range("S7").Value = 0.3448275862068966
# Row 8
$ws.Range("B8").Value = 0.09217877094972067
$ws.Range("D8").Value = 0.01675977653631285
$ws.Range("F8").Value = 0.05307262569832402
$ws.Range("J8").Value = 0.1508379888268156
$ws.Range("O8").Value = 0.01955307262569832
$ws.Range("Q8").Value = 0.1564245810055866
$ws.Range("R8").Value = 0.1005586592178771
$ws.Range("S8").Value = 0.4106145251396648
# Row 9
$ws.Range("B9").Value = 0.08074534161490683
$ws.Range("D9").Value = 0.0124223602484472
$ws.Range("F9").Value = 0.03105590062111801
$ws.Range("J9").Value = 0.1304347826086956
$ws.Range("O9").Value = 0.02484472049689441
$ws.Range("Q9").Value = 0.1925465838509317
$ws.Range("R9").Value = 0.1118012422360248
$ws.Range("S9").Value = 0.4161490683229814
# Row 10
$ws.Range("B10").Value = 0.08399646330680813
$ws.Range("D10").Value = 0.01768346595932803
$ws.Range("E10").Value = 0.0008841732979664014
$ws.Range("F10").Value = 0.07869142351900972
$ws.Range("J10").Value = 0.1202475685234306
$ws.Range("O10").Value = 0.02210433244916004
$ws.Range("Q10").Value = 0.2298850574712644
$ws.Range("R10").Value = 0.1043324491600354
$ws.Range("S10").Value = 0.3421750663129973
# Row 11
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.09863945578231292
$ws.Range("K11").Value = 0.1836734693877551
$ws.Range("L11").Value = 0.5680272108843537
$ws.Range("S11").Value = 0.006802721088435374
# Row 12
$ws.Range("G12").Value = 0.7247191011235955
$ws.Range("J12").Value = 0.1797752808988764
$ws.Range("K12").Value = 0.01685393258426966
$ws.Range("L12").Value = 0.05617977528089887
$ws.Range("S12").Value = 0.02247191011235955
# Row 13
$ws.Range("G13").Value = 0.48
$ws.Range("J13").Value = 0.4
$ws.Range("S13").Value = 0.12
# Row 15
$ws.Range("F15").Value = 0.0184331797235023
$ws.Range("H15").Value = 0.1244239631336406
$ws.Range("I15").Value = 0.05069124423963134
$ws.Range("J15").Value = 0.3640552995391705
$ws.Range("K15").Value = 0.06912442396313365
$ws.Range("M15").Value = 0.009216589861751152
$ws.Range("O15").Value = 0.06912442396313365
$ws.Range("S15").Value = 0.2949308755760369
# Row 16
$ws.Range("F16").Value = 0.03448275862068965
$ws.Range("H16").Value = 0.1586206896551724
$ws.Range("I16").Value = 0.1172413793103448
$ws.Range("J16").Value = 0.3931034482758621
$ws.Range("K16").Value = 0.1103448275862069
$ws.Range("M16").Value = 0.02758620689655172
$ws.Range("O16").Value = 0.04137931034482759
$ws.Range("S16").Value = 0.1172413793103448
# Row 17
$ws.Range("F17").Value = 0.01481481481481482
$ws.Range("H17").Value = 0.1901234567901235
$ws.Range("I17").Value = 0.08641975308641975
$ws.Range("J17").Value = 0.4271604938271605
$ws.Range("K17").Value = 0.09876543209876543
$ws.Range("M17").Value = 0.007407407407407408
$ws.Range("N17").Value = 0.002469135802469136
$ws.Range("O17").Value = 0.0691358024691358
$ws.Range("S17").Value = 0.1037037037037037
# Row 18
$ws.Range("F18").Value = 0.0196078431372549
$ws.Range("H18").Value = 0.196078431372549
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.3823529411764706
$ws.Range("K18").Value = 0.1372549019607843
$ws.Range("M18").Value = 0.00980392156862745
$ws.Range("O18").Value = 0.09313725490196079
$ws.Range("S18").Value = 0.0784313725490196
# Row 19
$ws.Range("F19").Value = 0.02171767028627838
$ws.Range("H19").Value = 0.1826258637709773
$ws.Range("I19").Value = 0.08094768015794669
$ws.Range("J19").Value = 0.3889437314906219
$ws.Range("K19").Value = 0.1303060217176703
$ws.Range("M19").Value = 0.01579466929911155
$ws.Range("N19").Value = 0.001974333662388944
$ws.Range("O19").Value = 0.07798617966436328
$ws.Range("S19").Value = 0.09970384995064166
